$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15; existing rows 15-36 shift down to 16-37.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new weekly price record.
$ws.Range("A15").Value = 8
$ws.Range("B15").Value = "Terminal La Palmera de La Serena"
$ws.Range("C15").Value = "Coquimbo"
$ws.Range("D15").Value = 44868
$ws.Range("E15").Value = 4
$ws.Range("F15").Value = "Fruta"
$ws.Range("G15").Value = 100101
$ws.Range("H15").Value = "Berries"
$ws.Range("I15").Value = 100101001
$ws.Range("J15").Value = "Arándano (blue)"
$ws.Range("K15").Value = "Sin especificar"
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 300
$ws.Range("N15").Value = 10000
$ws.Range("O15").Value = 11000
$ws.Range("P15").Value = 10500
$ws.Range("Q15").Value = "`$/bandeja 2 kilos"
$ws.Range("R15").Value = "Provincia de Limarí"
$ws.Range("S15").Value = 5250
$ws.Range("T15").Value = 2
